# Apply updated "想去人数" (want-to-go count) values across the
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsPerformance = $wb.Worksheets.Item("演出")
$wsAllTypes = $wb.Worksheets.Item("全部类型")

# 展览 (sheet1)
$wsExhibition.Range("F4").Value = 59
$wsExhibition.Range("F5").Value = 5042
$wsExhibition.Range("F7").Value = 87
$wsExhibition.Range("F8").Value = 299
$wsExhibition.Range("F9").Value = 52

# 演出 (sheet2)
$wsPerformance.Range("F2").Value = 127

# 全部类型 (sheet4)
$wsAllTypes.Range("F3").Value = 127
$wsAllTypes.Range("F8").Value = 59
$wsAllTypes.Range("F9").Value = 5042
$wsAllTypes.Range("F11").Value = 87
$wsAllTypes.Range("F13").Value = 299
$wsAllTypes.Range("F14").Value = 52
